$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# Column C: Riders
$ws.Range("C2").Value = 139
$ws.Range("C3").Value = 219
$ws.Range("C4").Value = 213
$ws.Range("C5").Value = 195
$ws.Range("C6").Value = 249
$ws.Range("C7").Value = 100
$ws.Range("C8").Value = 72

# Column D: Average
$ws.Range("D2").Value = 93.55
$ws.Range("D3").Value = 96.1
$ws.Range("D4").Value = 102.65
$ws.Range("D5").Value = 99.47
$ws.Range("D6").Value = 97.74
$ws.Range("D7").Value = 41.64
$ws.Range("D8").Value = 33.98
